# "add check in make changes"
#
# Someone added a new unavailability check ("NO") for Thursday-Evening
# (column F, row 3) on several people's personal constraint sheets, and
# then reshuffled the "shifts" master schedule so nobody is scheduled
# against their constraints anymore.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Personal constraint sheets: move/add the "NO" (unavailable) marks.
#    Each of these sheets uses row 2 = Morning, row 3 = Evening, and
#    columns B..H = Sunday..Saturday.
# ---------------------------------------------------------------------

$ws = $wb.Worksheets.Item("asaf")
$ws.Range("B2").Value = ""
$ws.Range("C3").Value = "NO"
$ws.Range("F3").Value = "NO"

$ws = $wb.Worksheets.Item("yoni")
$ws.Range("D3").Value = ""
$ws.Range("F3").Value = "NO"

$ws = $wb.Worksheets.Item("adir")
$ws.Range("C3").Value = ""
$ws.Range("F3").Value = "NO"

$ws = $wb.Worksheets.Item("stav")
$ws.Range("G2").Value = ""
$ws.Range("F3").Value = "NO"

$ws = $wb.Worksheets.Item("rotem")
$ws.Range("G2").Value = ""
$ws.Range("F3").Value = "NO"

# tair, michal and emilia keep their existing constraints unchanged.

# ---------------------------------------------------------------------
# 2) Master "shifts" schedule: re-assign workers so the updated
#    constraints above are respected.
# ---------------------------------------------------------------------

$ws = $wb.Worksheets.Item("shifts")

# Morning shift, first rotation (row 2)
$ws.Range("B2").Value = "tair"
$ws.Range("C2").Value = "adir"
$ws.Range("E2").Value = "yoni"
$ws.Range("F2").Value = "stav"

# Morning shift, second rotation (row 3)
$ws.Range("C3").Value = "yoni"
$ws.Range("D3").Value = "rotem"
$ws.Range("E3").Value = "rotem"
$ws.Range("F3").Value = "adir"
$ws.Range("G3").Value = "adir"

# Evening shift, first rotation (row 5)
$ws.Range("B5").Value = "adir"
$ws.Range("C5").Value = "adir"
$ws.Range("D5").Value = "rotem"
$ws.Range("E5").Value = "adir"
$ws.Range("F5").Value = "stav"
$ws.Range("H5").Value = "adir"

# Evening shift, second rotation (row 6)
$ws.Range("B6").Value = "asaf"
$ws.Range("C6").Value = "tair"
$ws.Range("D6").Value = "yoni"
$ws.Range("E6").Value = "tair"
$ws.Range("H6").Value = "asaf"
